$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.922.67'
$ws.Range('E2').Value = '  -0.28%  '
$ws.Range('D3').Value = '2.549.11'
$ws.Range('E3').Value = '  +0.03%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '304.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.49%  '
$ws.Range('E6').Value = '  +5.73%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.578'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.66%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -0.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.67'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0823'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.06%  '
$ws.Range('E12').Value = '  +4.44%  '
$ws.Range('E13').Value = '  -0.46%  '
$ws.Range('D14').Value = '2.941.19'
$ws.Range('E14').Value = '  +0.19%  '
$ws.Range('D15').Value = '2.567.36'
$ws.Range('E15').Value = '  +2.31%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.15'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +7.34%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.872'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.26%  '
$ws.Range('D18').Value = '42.951.08'
$ws.Range('E18').Value = '  -0.20%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.78'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.21%  '
$ws.Range('E20').Value = '  +1.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.58'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.90'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '254.38'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.01%  '
$ws.Range('E24').Value = '  +1.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.08'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '28.12'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.55%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.27'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '37.83'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.98%  '
$ws.Range('E30').Value = '  -1.78%  '
$ws.Range('E31').Value = '  +3.77%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '158.56'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.58'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +14.97%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.15'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0803'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.22%  '
$ws.Range('E36').Value = '  -2.21%  '
$ws.Range('E37').Value = '  -4.50%  '
$ws.Range('E38').Value = '  +1.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '25.64'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +9.51%  '
$ws.Range('E40').Value = '  -0.24%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.12'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +32.90%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.91'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('E43').Value = '  -0.79%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0307'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.33%  '
$ws.Range('D45').Value = '2.091.81'
$ws.Range('E45').Value = '  +0.55%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.998'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '86.55'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.32%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.97'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.82%  '
$ws.Range('B49').Value = 'ordi'
$ws.Range('C49').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '75.43'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +9.27%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.798.87'
$ws.Range('E50').Value = '  +0.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '103.34'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.42%  '
